$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45959
$ws.Range("B2").Value = 88.77
$ws.Range("C2").Value = 84.31999999999999
$ws.Range("D2").Value = 80.03
$ws.Range("E2").Value = 76.5
$ws.Range("F2").Value = 78.56999999999999
$ws.Range("G2").Value = 82.34
$ws.Range("H2").Value = 90.26000000000001
$ws.Range("I2").Value = 108.33
$ws.Range("J2").Value = 111.7
$ws.Range("K2").Value = 102.52
$ws.Range("L2").Value = 89.06
$ws.Range("M2").Value = 83.14
$ws.Range("N2").Value = 82.81999999999999
$ws.Range("O2").Value = 85.12
$ws.Range("P2").Value = 94.70999999999999
$ws.Range("Q2").Value = 103.02
$ws.Range("R2").Value = 108.51
$ws.Range("S2").Value = 123.55
$ws.Range("T2").Value = 137.32
$ws.Range("U2").Value = 152.38
$ws.Range("V2").Value = 159.89
$ws.Range("W2").Value = 147.6
$ws.Range("X2").Value = 126
$ws.Range("Y2").Value = 117.38
$ws.Range("Z2").Value = 104.74
$ws.Range("AB2").Value = 137.72
$ws.Range("AD2").Value = 153.74
$ws.Range("AF2").Value = 144.85
$ws.Range("AG2").Value = "0h-15h"
